# Sprint 44 test case report update
# - Bumps "Spint( 43)" labels to "Spint( 44)" across the per-day summary headers
# - Fills in the Day 2/3/4 execution/review totals
# - Moves the active selection/scroll position to reflect where work left off

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename all "Spint( 43) - Day N - Test Case Summary" headers to Sprint 44 ---
# The headers live in column B at the start of each day's block.
$headerRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 57)
for ($i = 0; $i -lt $headerRows.Length; $i++) {
    $day = $i + 1
    $row = $headerRows[$i]
    $ws.Range("B$row").Value = "Spint( 44) - Day $day - Test Case Summary"
}

# --- Fill in the newly completed totals for Day 2, Day 3 and Day 4 blocks ---
# Day 2 block (header row 8, data rows 9-11)
$ws.Range("C9").Value = 7075
$ws.Range("C10").Value = 2610
$ws.Range("C11").Value = 2610

# Day 3 block (header row 14, data rows 15-17)
$ws.Range("C15").Value = 7075
$ws.Range("C16").Value = 2610
$ws.Range("C17").Value = 2610

# Day 4 block (header row 20, data rows 21-23)
$ws.Range("C21").Value = 7075
$ws.Range("C22").Value = 2660
$ws.Range("C23").Value = 2660

# --- Update the view: move the active selection to C23 ---
$ws.Range("C23").Select() | Out-Null
